# "test data updated second time"
#
# The TestData sheet holds several small "customer id" lookup tables that
# reuse the same shared strings. This refresh swaps in a new batch of
# Stripe-style customer ids (cus_...), clears out one now-unused row, and
# moves the last block of ids to point at brand-new strings, leaving the
# sheet selection parked on the block that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# First lookup block (rows 14-18): refresh the four ids that are reused
# elsewhere on the sheet, then blank out the fifth (no longer needed).
$ws.Cells.Item(14, 1).Value = "cus_KHPAJc1Gvehntk"
$ws.Cells.Item(15, 1).Value = "cus_KHOlsXxLKlX9nS"
$ws.Cells.Item(16, 1).Value = "cus_KHOldthDPbluyt"
$ws.Cells.Item(17, 1).Value = "cus_KHOl86QmF04Gkg"
$ws.Cells.Item(18, 1).ClearContents()

# Second lookup block (rows 22-25) mirrors the same four ids - update in
# lockstep so both blocks keep showing the same values.
$ws.Cells.Item(22, 1).Value = "cus_KHPAJc1Gvehntk"
$ws.Cells.Item(23, 1).Value = "cus_KHOlsXxLKlX9nS"
$ws.Cells.Item(24, 1).Value = "cus_KHOldthDPbluyt"
$ws.Cells.Item(25, 1).Value = "cus_KHOl86QmF04Gkg"

# Third lookup block (rows 29-31): brand-new ids not used anywhere else.
$ws.Cells.Item(29, 1).Value = "cus_KHOfAbofEB8AAm"
$ws.Cells.Item(30, 1).Value = "cus_KHOf8fZSlUWWzy"
$ws.Cells.Item(31, 1).Value = "cus_KHOfRPNKXnm5rP"

# Leave the selection on the block that was just refreshed.
[void]$ws.Range("A29:A31").Select()
